# feat: add 2022-Q1 data
#
# Before:  Sheets = "2021-Q4", "总计"
# After:   Sheets = "2021-Q4", "2022-Q1", "总计"
#   - new "2022-Q1" sheet gets the same A1:H5 fund-holding table shape as "2021-Q4"
#   - "总计" sheet gets a new row (2022-Q1 / 4 / 1.23) inserted above the existing
#     2021-Q4 summary row

$wb = $excel.ActiveWorkbook

$ws2021 = $wb.Worksheets.Item("2021-Q4")

# NOTE: do NOT cache a reference to the "总计" sheet yet -- Worksheets.Item(...)
# references are positional, and Worksheets.Add() below shifts every sheet
# that comes after the insertion point. We re-fetch "总计" by name AFTER the
# new sheet has been inserted (see step 4).

# ---------------------------------------------------------------------------
# Helper: write a value into a cell and FORCE it to be stored as literal text
# (keeps leading zeros / decimal-looking strings like "011230" or "17.18" as
# text instead of letting Excel auto-coerce them to numbers). We do this by
# writing a formula that evaluates to the literal string, then converting the
# whole (already-written) range to static values in one shot with
# PasteSpecial(xlPasteValues) -- this does NOT touch cell formatting/styles.
# ---------------------------------------------------------------------------
function Set-LiteralTextFormula($range, [string]$value) {
    $escaped = $value -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
}

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right after "2021-Q4" (and therefore
#    right before "总计").
# ---------------------------------------------------------------------------
$ws2022 = $wb.Worksheets.Add($null, $ws2021)
$ws2022.Name = "2022-Q1"

# Clone the exact layout/formatting of "2021-Q4" onto the new sheet -- this
# reproduces the header style + index-column style without creating any new
# style records. Column A's header cell (A1) is genuinely empty on
# "2021-Q4", so copy it separately from the index values (A2:A5) to avoid
# materialising a spurious empty <c r="A1"/> element.
$ws2021.Range("B1:H5").Copy($ws2022.Range("B1"))
$ws2021.Range("A2:A5").Copy($ws2022.Range("A2"))

# ---------------------------------------------------------------------------
# 2) Overwrite the header row text.
# ---------------------------------------------------------------------------
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3) Overwrite the data rows.
#    Columns A & H are genuine numbers; B-G must stay as literal text
#    (fund codes have leading zeros, and the numeric-looking figures in the
#    source data are text cells, matching column D-G on the "2021-Q4" sheet).
# ---------------------------------------------------------------------------
$rows = @(
    @{ A = 0; B = "011230"; C = "创金合信数字经济主题股票C";             D = "17.18"; E = "92.17"; F = "4.07"; G = "0.6992"; H = 6 },
    @{ A = 1; B = "011229"; C = "创金合信数字经济主题股票A";             D = "12.18"; E = "92.17"; F = "4.07"; G = "0.4957"; H = 6 },
    @{ A = 2; B = "970020"; C = "信达价值精选一年持有期灵活配置混合A"; D = "0.64";  E = "56.02"; F = "3.14"; G = "0.0201"; H = 8 },
    @{ A = 3; B = "970021"; C = "信达价值精选一年持有期灵活配置混合B"; D = "0.53";  E = "56.02"; F = "3.14"; G = "0.0166"; H = 8 }
)

$r = 2
foreach ($row in $rows) {
    $ws2022.Cells.Item($r, 1).Value = $row.A
    Set-LiteralTextFormula $ws2022.Cells.Item($r, 2) $row.B
    $ws2022.Cells.Item($r, 3).Value = $row.C
    Set-LiteralTextFormula $ws2022.Cells.Item($r, 4) $row.D
    Set-LiteralTextFormula $ws2022.Cells.Item($r, 5) $row.E
    Set-LiteralTextFormula $ws2022.Cells.Item($r, 6) $row.F
    Set-LiteralTextFormula $ws2022.Cells.Item($r, 7) $row.G
    $ws2022.Cells.Item($r, 8).Value = $row.H
    $r++
}

# Convert the text-trick formulas (column B, D:G) into static literal values
# in one bulk operation -- preserves styles, removes the formulas.
$textRange = $ws2022.Range("B2:G5")
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues

# ---------------------------------------------------------------------------
# 4) Update the "总计" sheet: the existing single data row ("2021-Q4") moves
#    down to row 3, and a new row 2 is written for "2022-Q1".
#    Re-fetch the sheet by name now that the sheet collection has settled
#    (Worksheets.Add() above shifts positional references taken earlier).
#    We avoid Rows.Insert() here because it pulls in formatting from the
#    row above and creates a spurious new style record; instead the row 3
#    index cell's style is cloned directly from row 2's.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))

# Row 3 = old "2021-Q4" totals row, shifted down.
$wsTotal.Range("A3").Value = 1
Set-LiteralTextFormula $wsTotal.Range("B3") "2021-Q4"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.2

# Row 2 = new "2022-Q1" totals row.
$wsTotal.Range("A2").Value = 0
Set-LiteralTextFormula $wsTotal.Range("B2") "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 1.23

$totalTextRange = $wsTotal.Range("B2:B3")
$totalTextRange.Copy()
$totalTextRange.PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
